$p = $ppt.ActivePresentation
$p.ApplyTheme("C:\theme1.xml")
